$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Row 2
$ws.Cells.Item(2, 4).Value = '49.341.69'
$ws.Cells.Item(2, 5).Value = '  +2.71%  '

# Row 3
$ws.Cells.Item(3, 4).Value = '2.634.99'
$ws.Cells.Item(3, 5).Value = '  +5.77%  '

# Row 4
$ws.Cells.Item(4, 5).Value = '  +0.24%  '

# Row 5
$ws.Cells.Item(5, 4).Value = '111.17'
$ws.Cells.Item(5, 5).Value = '  +6.01%  '

# Row 6
$ws.Cells.Item(6, 4).Value = '323.03'
$ws.Cells.Item(6, 5).Value = '  +1.52%  '

# Row 7
$ws.Cells.Item(7, 4).Value = '0.521'
$ws.Cells.Item(7, 5).Value = '  +0.37%  '

# Row 8
$ws.Cells.Item(8, 4).Value = '1.00'
$ws.Cells.Item(8, 5).Value = '  +0.20%  '

# Row 9
$ws.Cells.Item(9, 4).Value = '0.541'
$ws.Cells.Item(9, 5).Value = '  +1.12%  '

# Row 10
$ws.Cells.Item(10, 4).Value = '39.77'
$ws.Cells.Item(10, 5).Value = '  +2.66%  '

# Row 11
$ws.Cells.Item(11, 4).Value = '19.88'
$ws.Cells.Item(11, 5).Value = '  -1.87%  '

# Row 12
$ws.Cells.Item(12, 4).Value = '0.0810'
$ws.Cells.Item(12, 5).Value = '  +1.31%  '

# Row 13
$ws.Cells.Item(13, 5).Value = '  -0.45%  '

# Row 14
$ws.Cells.Item(14, 4).Value = '7.22'
$ws.Cells.Item(14, 5).Value = '  +2.31%  '

# Row 15
$ws.Cells.Item(15, 4).Value = '3.063.65'
$ws.Cells.Item(15, 5).Value = '  +6.58%  '

# Row 16
$ws.Cells.Item(16, 4).Value = '2.662.04'
$ws.Cells.Item(16, 5).Value = '  +6.41%  '

# Row 17
$ws.Cells.Item(17, 4).Value = '0.860'
$ws.Cells.Item(17, 5).Value = '  +3.33%  '

# Row 18
$ws.Cells.Item(18, 4).Value = '49.409.15'
$ws.Cells.Item(18, 5).Value = '  +3.16%  '

# Row 19
$ws.Cells.Item(19, 4).Value = '12.84'
$ws.Cells.Item(19, 5).Value = '  +0.56%  '

# Row 20
$ws.Cells.Item(20, 4).Value = '6.65'
$ws.Cells.Item(20, 5).Value = '  +1.84%  '

# Row 21
$ws.Cells.Item(21, 4).Value = '2.89'
$ws.Cells.Item(21, 5).Value = '  -0.73%  '

# Row 22
$ws.Cells.Item(22, 4).Value = '0.0₃0941'
$ws.Cells.Item(22, 5).Value = '  +1.45%  '

# Row 23
$ws.Cells.Item(23, 4).Value = '271.68'
$ws.Cells.Item(23, 5).Value = '  -3.39%  '

# Row 24
$ws.Cells.Item(24, 4).Value = '70.16'
$ws.Cells.Item(24, 5).Value = '  -0.86%  '

# Row 25
$ws.Cells.Item(25, 4).Value = '2.53'
$ws.Cells.Item(25, 5).Value = '  +1.68%  '

# Row 26
$ws.Cells.Item(26, 4).Value = '26.27'
$ws.Cells.Item(26, 5).Value = '  +2.40%  '

# Row 27
$ws.Cells.Item(27, 4).Value = '0.998'
$ws.Cells.Item(27, 5).Value = '  -0.03%  '

# Row 28
$ws.Cells.Item(28, 4).Value = '10.07'
$ws.Cells.Item(28, 5).Value = '  +5.13%  '

# Row 29
$ws.Cells.Item(29, 5).Value = '  +0.95%  '

# Row 30
$ws.Cells.Item(30, 4).Value = '35.31'
$ws.Cells.Item(30, 5).Value = '  +2.25%  '

# Row 31
$ws.Cells.Item(31, 4).Value = '0.138'
$ws.Cells.Item(31, 5).Value = '  -0.92%  '

# Row 32
$ws.Cells.Item(32, 4).Value = '49.48'
$ws.Cells.Item(32, 5).Value = '  +0.92%  '

# Row 33
$ws.Cells.Item(33, 4).Value = '5.44'
$ws.Cells.Item(33, 5).Value = '  +3.68%  '

# Row 34
$ws.Cells.Item(34, 2).Value = 'FirstDigitalUSD'
$ws.Cells.Item(34, 3).Value = 'https://coinranking.com/coin/cpjRxjFYD+firstdigitalusd-fdusd'
$ws.Cells.Item(34, 4).Value = '1.00'
$ws.Cells.Item(34, 5).Value = '  -0.03%  '

# Row 35
$ws.Cells.Item(35, 2).Value = 'Celestia'
$ws.Cells.Item(35, 3).Value = 'https://coinranking.com/coin/YQcD0lBl7+celestia-tia'
$ws.Cells.Item(35, 4).Value = '19.11'
$ws.Cells.Item(35, 5).Value = '  -0.32%  '

# Row 36
$ws.Cells.Item(36, 4).Value = '0.0793'
$ws.Cells.Item(36, 5).Value = '  +3.25%  '

# Row 37
$ws.Cells.Item(37, 4).Value = '4.97'
$ws.Cells.Item(37, 5).Value = '  +10.78%  '

# Row 38
$ws.Cells.Item(38, 4).Value = '2.03'
$ws.Cells.Item(38, 5).Value = '  +4.49%  '

# Row 39
$ws.Cells.Item(39, 4).Value = '3.09'
$ws.Cells.Item(39, 5).Value = '  +7.58%  '

# Row 40
$ws.Cells.Item(40, 4).Value = '124.96'
$ws.Cells.Item(40, 5).Value = '  +4.25%  '

# Row 41
$ws.Cells.Item(41, 4).Value = '0.110'
$ws.Cells.Item(41, 5).Value = '  -0.31%  '

# Row 42
$ws.Cells.Item(42, 2).Value = 'WEMIXToken'
$ws.Cells.Item(42, 3).Value = 'https://coinranking.com/coin/08CsQa-Ov+wemixtoken-wemix'
$ws.Cells.Item(42, 4).Value = '2.21'
$ws.Cells.Item(42, 5).Value = '  +0.23%  '

# Row 43
$ws.Cells.Item(43, 2).Value = 'EnergySwap'
$ws.Cells.Item(43, 3).Value = 'https://coinranking.com/coin/SbWqqTui-+energyswap-ens'
$ws.Cells.Item(43, 4).Value = '21.97'
$ws.Cells.Item(43, 5).Value = '  +1.69%  '

# Row 44
$ws.Cells.Item(44, 4).Value = '0.0313'
$ws.Cells.Item(44, 5).Value = '  +4.96%  '

# Row 45
$ws.Cells.Item(45, 4).Value = '2.094.34'
$ws.Cells.Item(45, 5).Value = '  +5.42%  '

# Row 46
$ws.Cells.Item(46, 4).Value = '3.23'
$ws.Cells.Item(46, 5).Value = '  +4.02%  '

# Row 47
$ws.Cells.Item(47, 5).Value = '  +10.79%  '

# Row 48
$ws.Cells.Item(48, 5).Value = '  +5.04%  '

# Row 49
$ws.Cells.Item(49, 4).Value = '8.90'
$ws.Cells.Item(49, 5).Value = '  -0.26%  '

# Row 50
$ws.Cells.Item(50, 4).Value = '5.23'
$ws.Cells.Item(50, 5).Value = '  +2.80%  '

# Row 51
$ws.Cells.Item(51, 4).Value = '58.42'
$ws.Cells.Item(51, 5).Value = '  +4.41%  '
